# Auto-generated Excel COM-interop edit script
# Updates cryptos list Price (D) and Volume(1h) (E) columns per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.671.24'
$ws.Range("E2").Value = '  +1.84%  '
$ws.Range("D3").Value = '1.634.35'
$ws.Range("E3").Value = '  +1.91%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''212.56'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("E6").Value = '  +1.91%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = '''0.251'
$ws.Range("E8").Value = '  +1.36%  '
$ws.Range("D9").Value = '''0.0623'
$ws.Range("D10").Value = '''18.98'
$ws.Range("E10").Value = '  +4.34%  '
$ws.Range("E11").Value = '  +2.72%  '
$ws.Range("D12").Value = '1.863.52'
$ws.Range("E12").Value = '  +1.97%  '
$ws.Range("D13").Value = '1.644.36'
$ws.Range("E13").Value = '  +2.62%  '
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("D15").Value = '''0.524'
$ws.Range("E15").Value = '  +2.65%  '
$ws.Range("D16").Value = '26.675.82'
$ws.Range("E16").Value = '  +1.93%  '
$ws.Range("D17").Value = '''62.94'
$ws.Range("E17").Value = '  +1.93%  '
$ws.Range("D18").Value = '0.0₃0739'
$ws.Range("E18").Value = '  +1.89%  '
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").Value = '''208.05'
$ws.Range("E20").Value = '  +3.90%  '
$ws.Range("D21").Value = '''4.30'
$ws.Range("E21").Value = '  +0.82%  '
$ws.Range("D22").Value = '''9.38'
$ws.Range("E22").Value = '  +1.24%  '
$ws.Range("D23").Value = '''6.15'
$ws.Range("E23").Value = '  +2.71%  '
$ws.Range("E24").Value = '  +2.56%  '
$ws.Range("D25").Value = '''146.34'
$ws.Range("E25").Value = '  +1.51%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("E27").Value = '  -0.82%  '
$ws.Range("D28").Value = '''6.74'
$ws.Range("E28").Value = '  +2.96%  '
$ws.Range("D29").Value = '''15.37'
$ws.Range("E29").Value = '  +1.35%  '
$ws.Range("E30").Value = '  +5.60%  '
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("E32").Value = '  +1.05%  '
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("E34").Value = '  +1.23%  '
$ws.Range("D35").Value = '''1.49'
$ws.Range("E35").Value = '  +0.94%  '
$ws.Range("D36").Value = '1.166.29'
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("E37").Value = '  -0.41%  '
$ws.Range("D38").Value = '''0.807'
$ws.Range("E38").Value = '  +2.92%  '
$ws.Range("E40").Value = '  +1.74%  '
$ws.Range("E41").Value = '  +0.15%  '
$ws.Range("E42").Value = '  +1.63%  '
$ws.Range("E43").Value = '  +1.58%  '
$ws.Range("D44").Value = '1.773.99'
$ws.Range("E44").Value = '  +1.93%  '
$ws.Range("D45").Value = '''92.38'
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("E46").Value = '  +2.32%  '
$ws.Range("E47").Value = '  +7.86%  '
$ws.Range("D48").Value = '''54.78'
$ws.Range("E48").Value = '  +1.53%  '
$ws.Range("E49").Value = '  +1.42%  '
$ws.Range("D50").Value = '''0.409'
$ws.Range("E50").Value = '  +0.50%  '
$ws.Range("E51").Value = '  +4.64%  '
